# Apply "Changes in Shift Tag" edits to the Shifts worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update shift-tag text content for the affected cells
$ws.Range("B3").Value2 = "A`nNurse`nR  1.0`nA 0.0`nB`nNurse`nR  1.0`nA 0.0"
$ws.Range("D3").Value2 = "C`nNurse`nR  1.0`nA 0.0"
$ws.Range("E3").Value2 = "C`nNurse`nR  1.0`nA 0.0"
$ws.Range("B5").Value2 = "A/B Wing`nCNA`nR  3.0`nA 2.0`nC/D Wing`nCNA`nR  2.0`nA 1.0`nPrimrose`nCNA`nR  1.0`nA 0.0"
$ws.Range("C5").Value2 = "A/B Wing`nCNA`nR  3.0`nA 2.5`nC/D Wing`nCNA`nR  2.0`nA 1.0`nPrimrose`nCNA`nR  1.0`nA 0.0"
$ws.Range("G5").Value2 = "A/B Wing`nCNA`nR  3.0`nA 2.0`nC/D Wing`nCNA`nR  2.0`nA 1.0`nPrimrose`nCNA`nR  1.0`nA 0.0"
$ws.Range("H5").Value2 = "A/B Wing`nCNA`nR  3.0`nA 2.0`nC/D Wing`nCNA`nR  2.0`nA 1.0`nPrimrose`nCNA`nR  1.0`nA 0.0"
$ws.Range("C7").Value2 = "C/D Wing`nCNA`nR  1.0`nA 0.0"
$ws.Range("D7").Value2 = "C/D Wing`nCNA`nR  1.0`nA 0.0`nJCR`nCNA`nR  1.0`nA 0.0`nPrimrose`nCNA`nR  1.0`nA 0.0"
$ws.Range("E7").Value2 = "C/D Wing`nCNA`nR  1.0`nA 0.0`nJCR`nCNA`nR  1.0`nA 0.0`nPrimrose`nCNA`nR  1.0`nA 0.0"
$ws.Range("F7").Value2 = "C/D Wing`nCNA`nR  1.0`nA 0.0`nJCR`nCNA`nR  1.0`nA 0.0"

# These cells were (and remain) intentionally blank; re-clear them so the
# load/save round-trip does not leave a stray shared-string reference behind.
$ws.Range("D4").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("G6").ClearContents()

Write-Output "Done applying shift tag changes"
